# Fill in the "Absent" column (H) for the consolidated attendance report.
# Rows 5, 6, 9 and 13 were missing/incorrect Absent values; set them
# explicitly to complete the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H13").Value = 0
